$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$rng = $ws.Range("A$row`:M$row")
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "fatima3"
$ws.Cells.Item($row, 2).Value  = "06@gmail.com"
$ws.Cells.Item($row, 3).Value  = "9234556789"
$ws.Cells.Item($row, 4).Value  = "Intermediate"
$ws.Cells.Item($row, 5).Value  = "2024"
$ws.Cells.Item($row, 6).Value  = "Australia"
$ws.Cells.Item($row, 7).Value  = "diploma"
$ws.Cells.Item($row, 8).Value  = "nj"
$ws.Cells.Item($row, 9).Value  = "No"
$ws.Cells.Item($row, 10).Value = "Yes"
$ws.Cells.Item($row, 11).Value = "http://localhost:5173/services/applyfordiplomacourses"
$ws.Cells.Item($row, 12).Value = "30/10/2025, 2:39:14 pm"
$ws.Cells.Item($row, 13).Value = "69032b3a322bf809ad24b953"
